$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the new "basketData" sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "formData"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "basketData"

# --- Populate basketData ---
$ws2.Range("A1").Value = "testCase1"
$ws2.Range("B1").Value = "Product"
$ws2.Range("C1").Value = "Count"

$ws2.Range("B2").Value = "Okulary"
$ws2.Range("C2").Value = 2
$ws2.Range("B3").Value = "Kabel"
$ws2.Range("C3").Value = 4
$ws2.Range("B4").Value = "Słuchawki"
$ws2.Range("C4").Value = 1
$ws2.Range("B5").Value = "Piłka"
$ws2.Range("C5").Value = 5

$ws2.Range("D6").Value = "testCase1"

$ws2.Range("A8").Value = "testCase2"
$ws2.Range("B8").Value = "Product"
$ws2.Range("C8").Value = "Count"

$ws2.Range("B9").Value = "Aparat"
$ws2.Range("C9").Value = 2
$ws2.Range("B10").Value = "Zeszyt"
$ws2.Range("C10").Value = 4
$ws2.Range("B11").Value = "Kostka"
$ws2.Range("C11").Value = 1

$ws2.Range("D12").Value = "testCase2"

# The "Count" column reuses the workbook's existing blue-font cell style
# (the same one applied to the hyperlinked e-mails on formData).
$ws2.Range("C2:C5").Font.Color = 16711680
$ws2.Range("C9:C11").Font.Color = 16711680
$ws2.Range("C12").Font.Color = 16711680

# Trailing formatted-but-empty cells that round out the original layout
# (touching NumberFormat registers the cell without fabricating a new style).
$ws2.Range("D2:D5").NumberFormat = "General"
$ws2.Range("D9:D11").NumberFormat = "General"
$ws2.Range("B12").NumberFormat = "General"
$ws2.Range("D13").NumberFormat = "General"

# --- Column widths ---
# formData: only column A nudges from 21.54 -> ~21.56 (closest the engine's
# internal 1/6-character width grid can reach).
$ws1.Columns.Item(1).ColumnWidth = 20.666666666666668

# basketData
$ws2.Columns.Item(2).ColumnWidth = 24.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 26.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 23.666666666666668

# --- Selection states to match target workbook ---
$ws1.Range("E6").Select()
$ws2.Range("B15").Select()
